$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state (after the GitHub Actions crypto-price refresh) for columns B (Coin),
# C (Link), D (Price) and E (Volume(1h)) across data rows 2-51.
# A new coin (LEO) was inserted at row 23, shifting the following coins down by one
# row; the former last row (Cronos) drops off the bottom of the fixed-size table.
$coinRows = @(
    @{Row = 2; B = "Bitcoin"; C = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D = "66.791.71"; E = "  -2.38%  "}
    @{Row = 3; B = "Ethereum"; C = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D = "2.453.36"; E = "  -3.61%  "}
    @{Row = 4; B = "TetherUSD"; C = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D = "1.00"; E = "  -0.06%  "}
    @{Row = 5; B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "580.45"; E = "  -2.16%  "}
    @{Row = 6; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "167.65"; E = "  -4.71%  "}
    @{Row = 7; B = "USDC"; C = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D = "1.00"; E = "  +0.01%  "}
    @{Row = 8; B = "XRP"; C = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D = "0.509"; E = "  -3.21%  "}
    @{Row = 9; B = "LidoStakedEther"; C = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"; D = "2.453.06"; E = "  -3.70%  "}
    @{Row = 10; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "0.133"; E = "  -3.78%  "}
    @{Row = 11; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "0.165"; E = "  -1.00%  "}
    @{Row = 12; B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "4.87"; E = "  -3.36%  "}
    @{Row = 13; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "0.326"; E = "  -5.76%  "}
    @{Row = 14; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "2.903.05"; E = "  -1.82%  "}
    @{Row = 15; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "25.14"; E = "  -5.60%  "}
    @{Row = 16; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "66.794.20"; E = "  -2.46%  "}
    @{Row = 17; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.0000167"; E = "  -5.72%  "}
    @{Row = 18; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "2.413.81"; E = "  -7.55%  "}
    @{Row = 19; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "10.91"; E = "  -8.54%  "}
    @{Row = 20; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "7.36"; E = "  -8.63%  "}
    @{Row = 21; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "348.73"; E = "  -5.94%  "}
    @{Row = 22; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "4.00"; E = "  -4.20%  "}
    @{Row = 23; B = "LEO"; C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D = "5.94"; E = "  -2.17%  "}
    @{Row = 24; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "0.998"; E = "  -0.21%  "}
    @{Row = 25; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "68.53"; E = "  -4.72%  "}
    @{Row = 26; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "4.18"; E = "  -8.81%  "}
    @{Row = 27; B = "SuiNetwork"; C = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"; D = "1.79"; E = "  -6.34%  "}
    @{Row = 28; B = "Aptos"; C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D = "9.03"; E = "  -9.39%  "}
    @{Row = 29; B = "Binance-PegBSC-USD"; C = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"; D = "1.00"; E = "  -39.29%  "}
    @{Row = 30; B = "WrappedeETH"; C = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"; D = "2.584.83"; E = "  -3.34%  "}
    @{Row = 31; B = "PEPE"; C = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D = "0.0₃0888"; E = "  -8.62%  "}
    @{Row = 32; B = "Bittensor"; C = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D = "504.02"; E = "  -6.11%  "}
    @{Row = 33; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "7.57"; E = "  -9.03%  "}
    @{Row = 34; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "1.75"; E = "  -6.50%  "}
    @{Row = 35; B = "Fetch.AI"; C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D = "1.21"; E = "  -8.62%  "}
    @{Row = 36; B = "FirstDigitalUSD"; C = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D = "1.00"; E = "  -0.09%  "}
    @{Row = 37; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "158.13"; E = "  -1.43%  "}
    @{Row = 38; B = "Kaspa"; C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D = "0.113"; E = "  -12.80%  "}
    @{Row = 39; B = "WhiteBITCoin"; C = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"; D = "18.64"; E = "  +0.01%  "}
    @{Row = 40; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "18.13"; E = "  -6.02%  "}
    @{Row = 41; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "1.31"; E = "  -9.59%  "}
    @{Row = 42; B = "USDe"; C = "https://coinranking.com/coin/exbfr2U-0+usde-usde"; D = "1.00"; E = "  -0.40%  "}
    @{Row = 43; B = "Stacks"; C = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"; D = "1.66"; E = "  -7.48%  "}
    @{Row = 44; B = "RenderToken"; C = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D = "4.74"; E = "  -8.42%  "}
    @{Row = 45; B = "PolygonEcosystemToken"; C = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"; D = "0.322"; E = "  -8.14%  "}
    @{Row = 46; B = "dogwifhat"; C = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D = "2.34"; E = "  -8.01%  "}
    @{Row = 47; B = "OKB"; C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D = "38.54"; E = "  -2.31%  "}
    @{Row = 48; B = "Aave"; C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D = "139.71"; E = "  -6.14%  "}
    @{Row = 49; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "3.40"; E = "  -8.76%  "}
    @{Row = 50; B = "ARBITRUM"; C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D = "0.503"; E = "  -9.27%  "}
    @{Row = 51; B = "BabyDogeCoin"; C = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"; D = "0.0₆0251"; E = "  -11.04%  "}
)

foreach ($coin in $coinRows) {
    $bCell = $ws.Cells.Item($coin.Row, 2)
    $cCell = $ws.Cells.Item($coin.Row, 3)
    $dCell = $ws.Cells.Item($coin.Row, 4)
    $eCell = $ws.Cells.Item($coin.Row, 5)

    # Force text format so numeric-looking strings (prices, percentages) are not
    # auto-converted to numbers by Excel.
    $bCell.NumberFormat = "@"
    $cCell.NumberFormat = "@"
    $dCell.NumberFormat = "@"
    $eCell.NumberFormat = "@"

    $bCell.Value = $coin.B
    $cCell.Value = $coin.C
    $dCell.Value = $coin.D
    $eCell.Value = $coin.E
}
